$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting existing rows 116:178 down to 117:179
$ws.Rows("116").Insert()

# Populate the new row 116 with the template values used throughout this data block,
# with the specific new data values for D, J, K, L, M, P
$ws.Cells.Item(116, 1).Value = 8
$ws.Cells.Item(116, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(116, 3).Value = "Coquimbo"
$ws.Cells.Item(116, 4).Value = 44873
$ws.Cells.Item(116, 5).Value = 4
$ws.Cells.Item(116, 6).Value = 100112044
$ws.Cells.Item(116, 7).Value = "Perejil"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 2400
$ws.Cells.Item(116, 11).Value = 1500
$ws.Cells.Item(116, 12).Value = 2000
$ws.Cells.Item(116, 13).Value = 1750
$ws.Cells.Item(116, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(116, 16).Value = 1167
$ws.Cells.Item(116, 17).Value = 1.5
$ws.Cells.Item(116, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same number format / style as the rest of column D
$ws.Cells.Item(116, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat
